# LOM3205.xlsx update: add a dedicated "Docentes responsaveis" (responsible
# professors) block, give "Objetivos/Objectives" their own descriptive text
# (previously the professor names had leaked into those cells), and refresh
# the syllabus/method/criteria/bibliography text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 5 rows for the new "Docentes responsaveis:" block (1 label
# row + 4 professor rows) right after the Objectives rows.
$ws.Rows("12:16").Insert()

# Row 10
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Estudo formal da teoria dos campos eletromagnéticos independentes do tempo ou para situações quase-estáticas. Teoria das ondas eletromagnéticas."
$ws.Range("C10").Value = "Estudo formal da teoria dos campos eletromagnéticos independentes do tempo ou para situações quase-estáticas. Teoria das ondas eletromagnéticas."
$ws.Rows(10).RowHeight = 60

# Row 11
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."
$ws.Range("C11").Value = "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."
$ws.Rows(11).RowHeight = 60

# Row 12
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()

# Row 13
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C13").Value = "5840726 - Cristina Bormio Nunes"

# Row 14
$ws.Range("A14").Clear()
$ws.Range("B14").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C14").Value = "6495737 - Durval Rodrigues Junior"

# Row 15
$ws.Range("A15").Clear()
$ws.Range("B15").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C15").Value = "1341653 - Maria José Ramos Sandim"

# Row 16
$ws.Range("A16").Clear()
$ws.Range("B16").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("C16").Value = "1643715 - Paulo Atsushi Suzuki"

# Row 17
$ws.Range("A17").Value = "Programa resumido:"
$ws.Range("B17").Value = "Eletrostática. Magnetostática. Campos variantes no tempo. Equações de Maxwell. Ondas eletromagnéticas."
$ws.Range("C17").Value = "Eletrostática. Magnetostática. Campos variantes no tempo. Equações de Maxwell. Ondas eletromagnéticas."
$ws.Rows(17).RowHeight = 60

# Row 18
$ws.Range("A18").Value = "Short syllabus:"
$ws.Range("B18").Value = "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"
$ws.Range("C18").Value = "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"
$ws.Rows(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Programa:"
$ws.Range("B19").Value = "Eletrostática (campo eletrostático; potencial elétrico; trabalho e energia em eletrostática).  Técnicas especiais para a resolução da equação de Laplace (método das imagens; separação de variáveis). Campo elétrico da matéria (polarização elétrica; campo de objeto polarizado; cargas ligadas; deslocamento elétrico; dielétricos (lineares). Magnetostática (Lei de Lorentz; Lei de Biot-Savart; Lei de Ampére; vetor potencial magnético).  Campo magnético na matéria (magnetização; campos de objeto magnetizado; campo auxiliar H; Eletrodinâmica (força eletromotriz; indução eletromagnética; equações de Maxwell; lei de conservação de carga). Ondas eletromagnéticas (propagação no vácuo e na matéria; reflexão e transmissão), equação de ondas (planas)  e condições de contorno (interfaces). Radiação de dipolo elétrico."
$ws.Range("C19").Value = "Eletrostática (campo eletrostático; potencial elétrico; trabalho e energia em eletrostática).  Técnicas especiais para a resolução da equação de Laplace (método das imagens; separação de variáveis). Campo elétrico da matéria (polarização elétrica; campo de objeto polarizado; cargas ligadas; deslocamento elétrico; dielétricos (lineares). Magnetostática (Lei de Lorentz; Lei de Biot-Savart; Lei de Ampére; vetor potencial magnético).  Campo magnético na matéria (magnetização; campos de objeto magnetizado; campo auxiliar H; Eletrodinâmica (força eletromotriz; indução eletromagnética; equações de Maxwell; lei de conservação de carga). Ondas eletromagnéticas (propagação no vácuo e na matéria; reflexão e transmissão), equação de ondas (planas)  e condições de contorno (interfaces). Radiação de dipolo elétrico."
$ws.Rows(19).RowHeight = 120

# Row 20
$ws.Range("A20").Value = "Syllabus:"
$ws.Range("B20").Value = "Electrostatics (electrostatic field; electric potential; work and energy in electrostatics). Special techniques for solving the Laplace’s equation (method of images; separation of variables). Electric field in matter (electric polarization; polarized object field; bound charges; electric displacement; dielectrics (linear). Magnetostatics (Lorentz's law; Biot-Savart's law; Ampere's law; vector magnetic potential). Magnetic field in matter (magnetization; field of a magnetized object; auxiliary field H). Electrodynamics (electromotive force; electromagnetic induction; Maxwell's equations; law of conservation of charge). Electromagnetic waves (propagation in vacuum and in matter; reflection and transmission), plane wave equation and boundary conditions (interfaces). Electric dipole radiation."
$ws.Range("C20").Value = "Electrostatics (electrostatic field; electric potential; work and energy in electrostatics). Special techniques for solving the Laplace’s equation (method of images; separation of variables). Electric field in matter (electric polarization; polarized object field; bound charges; electric displacement; dielectrics (linear). Magnetostatics (Lorentz's law; Biot-Savart's law; Ampere's law; vector magnetic potential). Magnetic field in matter (magnetization; field of a magnetized object; auxiliary field H). Electrodynamics (electromotive force; electromagnetic induction; Maxwell's equations; law of conservation of charge). Electromagnetic waves (propagation in vacuum and in matter; reflection and transmission), plane wave equation and boundary conditions (interfaces). Electric dipole radiation."
$ws.Rows(20).RowHeight = 120

# Row 21
$ws.Range("A21").Value = "Avaliação:"

# Row 22
$ws.Range("A22").Value = "Método:"
$ws.Range("B22").Value = "Aulas expositivas e  exercícios comentados"
$ws.Range("C22").Value = "Aulas expositivas e  exercícios comentados"
$ws.Rows(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = "Critério:"
$ws.Range("B23").Value = "Média final calculada pelas notas de 2 provas (P1 e P2), seguindo os pesos MF=(P1+2*P2)/3, ou seja, peso 1 para a P1 e peso 2 para a P2."
$ws.Range("C23").Value = "Média final calculada pelas notas de 2 provas (P1 e P2), seguindo os pesos MF=(P1+2*P2)/3, ou seja, peso 1 para a P1 e peso 2 para a P2."
$ws.Rows(23).RowHeight = 60

# Row 24
$ws.Range("A24").Value = "Norma de recuperação:"
$ws.Range("B24").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C24").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows(24).RowHeight = 60

# Row 25
$ws.Range("A25").Value = "Bibliografia:"
$ws.Range("B25").Value = "CHENG,DAVID K.Field and Wave Electromagnetics. Addison Weslwy Publishing Company.1989.`nSLATER, J.C.; FRANK, N.H. Electromagnetism. McGraw-Hill, New York, 1974.`nMARION, J.B. Classical Electromagnetic Radiation. Academic Press, New York, 1965.`nBOHN, E.V. Introduction to electromagnetic fields and waves. Addison Wesley, 1968. `nREITZ, J.R.; MILFORD, F.J. Foundations of eletromagnetic theory. Addison Wesley, Publishing, Co. 1970. GRIFFITHS, D.J. Introduction to Electrodynamics. Prentice Hall, New York. 1998. `nRAMO, WHINNERY E VAN DUZER, Fields and Waves in Communication Electronics, Wiley."
$ws.Range("C25").Value = "CHENG,DAVID K.Field and Wave Electromagnetics. Addison Weslwy Publishing Company.1989.`nSLATER, J.C.; FRANK, N.H. Electromagnetism. McGraw-Hill, New York, 1974.`nMARION, J.B. Classical Electromagnetic Radiation. Academic Press, New York, 1965.`nBOHN, E.V. Introduction to electromagnetic fields and waves. Addison Wesley, 1968. `nREITZ, J.R.; MILFORD, F.J. Foundations of eletromagnetic theory. Addison Wesley, Publishing, Co. 1970. GRIFFITHS, D.J. Introduction to Electrodynamics. Prentice Hall, New York. 1998. `nRAMO, WHINNERY E VAN DUZER, Fields and Waves in Communication Electronics, Wiley."
$ws.Rows(25).RowHeight = 120

# Row 26
$ws.Range("A26").Value = "Requisitos:"

# Row 27
$ws.Range("B27").Value = "LOB1052 -  Cálculo III  (Requisito)`n"
$ws.Range("C27").Value = "LOB1052 -  Cálculo III  (Requisito)`n"
$ws.Rows(27).RowHeight = 30

# Row 28
$ws.Range("B28").Value = "LOB1053 -  Física III  (Requisito)`n"
$ws.Range("C28").Value = "LOB1053 -  Física III  (Requisito)`n"
$ws.Rows(28).RowHeight = 30
